# labor-timekeeper export update for Phil_Henderson_2026-01-12.xlsx
# - revert admin dev default employee id
# - reseed weekly timesheet with new customers (table was previously seeded
#   with sample/placeholder hours; now cleared back to empty defaults)
# - remove the extra seeded OT row (overtime sample row) entirely
# - autosave values now reflect Reg: 40 / OT: 0 and zeroed rate/total

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Weekly Timesheet"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Weekly Timesheet")

# Row 2 - 2026-01-12
$ws1.Range("B2").Value = "Bottomley"
$ws1.Range("C2").Value = 8
$ws1.Range("E2").Value = 0
$ws1.Range("F2").Value = 0

# Row 3 - 2026-01-13
$ws1.Range("B3").Value = "Bahin"
$ws1.Range("E3").Value = 0
$ws1.Range("F3").Value = 0

# Row 4 - 2026-01-14
$ws1.Range("B4").Value = "Schauer"
$ws1.Range("E4").Value = 0
$ws1.Range("F4").Value = 0

# Row 5 - 2026-01-15
$ws1.Range("B5").Value = "Hendricks"
$ws1.Range("E5").Value = 0
$ws1.Range("F5").Value = 0

# Row 6 - 2026-01-16
$ws1.Range("B6").Value = "Tubergen"
$ws1.Range("C6").Value = 8
$ws1.Range("E6").Value = 0
$ws1.Range("F6").Value = 0

# SUBTOTAL row (row 9) - hours total 44 -> 40, note text, total 1380 -> 0
$ws1.Range("C9").Value = 40
$ws1.Range("D9").Value = "Reg: 40 / OT: 0"
$ws1.Range("F9").Value = 0

# HOURLY SUBTOTAL row (row 12) - total 1380 -> 0
$ws1.Range("F12").Value = 0

# GRAND TOTAL row (row 14) - total 1380 -> 0
$ws1.Range("F14").Value = 0

# Remove the seeded OT row entirely (row 7)
$ws1.Rows("7").Delete()

# ---------------------------------------------------------------------
# Sheet 2: "Jason Schema"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Jason Schema")

# Revert the admin dev default employee id everywhere it's seeded
$ws2.Range("B2").Value = "emp_y716hily"
$ws2.Range("B3").Value = "emp_y716hily"
$ws2.Range("B4").Value = "emp_y716hily"
$ws2.Range("B5").Value = "emp_y716hily"
$ws2.Range("B6").Value = "emp_y716hily"

# Row 2 - 2026-01-12
$ws2.Range("D2").Value = "Bottomley"
$ws2.Range("E2").Value = 8
$ws2.Range("F2").Value = 0
$ws2.Range("G2").Value = 0
$ws2.Range("I2").Value = ""

# Row 3 - 2026-01-13
$ws2.Range("D3").Value = "Bahin"
$ws2.Range("F3").Value = 0
$ws2.Range("G3").Value = 0
$ws2.Range("I3").Value = ""

# Row 4 - 2026-01-14
$ws2.Range("D4").Value = "Schauer"
$ws2.Range("F4").Value = 0
$ws2.Range("G4").Value = 0
$ws2.Range("I4").Value = ""

# Row 5 - 2026-01-15
$ws2.Range("D5").Value = "Hendricks"
$ws2.Range("F5").Value = 0
$ws2.Range("G5").Value = 0
$ws2.Range("I5").Value = ""

# Row 6 - 2026-01-16
$ws2.Range("D6").Value = "Tubergen"
$ws2.Range("E6").Value = 8
$ws2.Range("F6").Value = 0
$ws2.Range("G6").Value = 0
$ws2.Range("I6").Value = ""

# Remove the seeded OT row entirely (row 7)
$ws2.Rows("7").Delete()

$wb.Save()
